$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Countries without divisions/factories data yet: clear the "0" placeholders
# in the Фабрики (C) / ВоенЗаводы (D) columns and recolor them like
# "Готово, доработать до альфы" (fillId 4 / style index 10), matching the
# formatting already used for the legend swatch in P9.

$donorBlank = $ws.Range("P9")      # style index 10 (fillId 4, no border)

foreach ($row in 2,4,5,6,7) {
    $donorBlank.Copy() | Out-Null
    $ws.Range("C" + $row + ":D" + $row).PasteSpecial(-4122) | Out-Null
    $ws.Range("C" + $row + ":D" + $row).ClearContents() | Out-Null
}

# --- СЛАНДСКО (row 5) now has 20 divisions -> recolor like "Почти готово"
# (fillId 5 / style index 9, matching legend swatch P8) and set the value.
$donorDiv = $ws.Range("P8")        # style index 9 (fillId 5, no border)
$donorDiv.Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").Value = 20

# --- New country row: ЭЛОНГАЛ (row 8), matching the layout/format of the
# other country rows (style copied from row 7 which is fully filled in).
$ws.Range("A8").Value = "ЭЛОНГАЛ"

$ws.Range("B7").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").Value = 0

$donorBlank.Copy() | Out-Null
$ws.Range("C8:D8").PasteSpecial(-4122) | Out-Null

$ws.Range("E7:F7").Copy() | Out-Null
$ws.Range("E8:F8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").Value = 2

$ws.Range("G7:K7").Copy() | Out-Null
$ws.Range("G8:K8").PasteSpecial(-4122) | Out-Null

# --- last selected cell in the sheet when the file was saved
$ws.Range("D9").Select() | Out-Null
